{"js": "// After the \"1h bars\" training-results list, two more results were logged\n// (40 epochs and 50 epochs). Add them as new paragraphs right after the\n// paragraph reading \"78 6820\", and move the \"_GoBack\" last-edit bookmark to\n// the end of the newly typed text (Word re-marks the _GoBack bookmark at the\n// spot of the most recent edit).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"78 6820\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find paragraph '78 6820'\");\n}\n\nconst p1 = target.insertParagraph(\"40 - 76 6720\", \"After\");\nconst p2 = p1.insertParagraph(\"50 - 323 28403\", \"After\");\n\n// Relocate the \"_GoBack\" bookmark (tracks Word's last edit position) to the\n// text that was just typed.\ncontext.document.deleteBookmark(\"_GoBack\");\np2.getRange(\"Content\").insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# After the \"1h bars\" training-results list, two more results were logged\n# (40 epochs and 50 epochs). Add them as new paragraphs right after the\n# paragraph reading \"78 6820\", and move the \"_GoBack\" last-edit bookmark to\n# the end of the newly typed text (Word re-marks the _GoBack bookmark at the\n# spot of the most recent edit).\n$d = $word.ActiveDocument\n\n# Locate the target paragraph robustly by its exact text (paragraph Range.Text\n# always carries the trailing paragraph mark).\n$targetIndex = -1\n$i = 1\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -eq \"78 6820`r\") {\n        $targetIndex = $i\n    }\n    $i = $i + 1\n}\nif ($targetIndex -eq -1) {\n    throw \"Could not find paragraph '78 6820'\"\n}\n\n$target = $d.Paragraphs.Item($targetIndex)\n$target.Range.InsertParagraphAfter()\n$p1 = $d.Paragraphs.Item($targetIndex + 1)\n$p1.Range.Text = \"40 - 76 6720\"\n\n$p1.Range.InsertParagraphAfter()\n$p2 = $d.Paragraphs.Item($targetIndex + 2)\n$p2.Range.Text = \"50 - 323 28403\"\n\n# Relocate the \"_GoBack\" last-edit bookmark onto the text just typed.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n$bmRange = $d.Range($p2.Range.Start, $p2.Range.End - 1)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n"}
